$d = $word.ActiveDocument

# Word constants
$wdHeaderFooterPrimary   = 1
$wdHeaderFooterFirstPage = 2

$sec = $d.Sections.Item(1)

# --- Header (first-page header, header1.xml): BTec_Logo-Orange
#     name: image1.jpg -> image2.jpg ----------------------------------------
$hdrFirst = $sec.Headers.Item($wdHeaderFooterFirstPage)
$ishHdr = $hdrFirst.Range.InlineShapes.Item(1)
if ($ishHdr.AlternativeText -eq "BTec_Logo-Orange") {
    $ishHdr2 = $ishHdr.Range.InlineShapes.Item(1)
    $ishHdr2.Name = "image2.jpg"
}

# --- Footer (default footer, footer2.xml): Pearson logo
#     name: image2.png -> image1.png ----------------------------------------
$ftrDefault = $sec.Footers.Item($wdHeaderFooterPrimary)
$ishFtrD = $ftrDefault.Range.InlineShapes.Item(1)
if ($ishFtrD.AlternativeText -like "*PearsonLogo.png") {
    $ishFtrD2 = $ishFtrD.Range.InlineShapes.Item(1)
    $ishFtrD2.Name = "image1.png"
}

# --- Footer (first-page footer, footer1.xml): Pearson logo
#     name: image2.png -> image1.png ----------------------------------------
$ftrFirst = $sec.Footers.Item($wdHeaderFooterFirstPage)
$ishFtrF = $ftrFirst.Range.InlineShapes.Item(1)
if ($ishFtrF.AlternativeText -like "*PearsonLogo.png") {
    $ishFtrF2 = $ishFtrF.Range.InlineShapes.Item(1)
    $ishFtrF2.Name = "image1.png"
}
